# daily auto push: 2026-02-22 18:52 UTC
# Insert 3 new rows of data (2026/02/22 19:00, 2026/02/22 22:00, 2026/02/23 02:00)
# just before the existing "2026/12/29" block, shifting all subsequent rows down
# by 3 (old last row 3188 -> new last row 3191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 3147 onward down by three rows.
$ws.Rows("3147:3149").Insert()

# Column A holds date-like text (e.g. "2026/02/22"); force it to stay text
# (matching the rest of the column) instead of being auto-converted to a
# date serial number.
$ws.Range("A3147:A3149").NumberFormat = "@"

$ws.Range("A3147").Value = "2026/02/22"
$ws.Range("B3147").Value = "日"
$ws.Range("C3147").Value = 19
$ws.Range("D3147").Value = 63

$ws.Range("A3148").Value = "2026/02/22"
$ws.Range("B3148").Value = "日"
$ws.Range("C3148").Value = 22
$ws.Range("D3148").Value = 60

$ws.Range("A3149").Value = "2026/02/23"
$ws.Range("B3149").Value = "月"
$ws.Range("C3149").Value = 2
$ws.Range("D3149").Value = 53

# Drop the "text" number format so the cells end up with no explicit style,
# matching the plain unstyled cells used throughout the rest of the sheet.
$ws.Range("A3147:A3149").Style = "Normal"

Write-Host "Inserted rows 3147-3149 (2026/02/22-23 data); sheet now spans A1:D3191."
